# EWU Dashboard - add "Impact-based Weighting" columns to the Weighting sheet
# Inserts 5 new weighting-method columns (B:F) in front of the existing
# weighting matrix, fills their header labels, and copies number formats
# from the neighbouring (former first) data column so the new blank cells
# look like the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weighting")

# --- insert 5 new columns before column B -------------------------------
$ws.Range("B1:F1").EntireColumn.Insert()

# --- copy formatting (number format, fill, borders, font, alignment) ---
# header row (row 1) - copy from the (new) first header cell, G1
$ws.Range("G1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)   # xlPasteFormats

# data rows (rows 2-17) - copy from the (new) first data column, G2:G17
$ws.Range("G2:G17").Copy()
$ws.Range("B2:F17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# --- match the column widths of the rest of the matrix ------------------
$ws.Range("B1:F1").EntireColumn.ColumnWidth = $ws.Range("G1").EntireColumn.ColumnWidth

# --- header labels for the new impact-based weighting columns -----------
$ws.Range("B1").Value2 = "Impact-based Weighting"
$ws.Range("C1").Value2 = "Impact-based Weighting + Importance and Robustness (Sala et al. 2018)"
$ws.Range("D1").Value2 = "Impact-based Weighting + Distance-to-Target (Castellani et al. 2016)"
$ws.Range("E1").Value2 = "Impact-based Weighting + Planetary Boundaries (Bjorn & Hauschild 2015)"
$ws.Range("F1").Value2 = "Impact-based Weighting + Quality and Maturity (ILCD 2011)"

# header row grows taller to fit the longer wrapped labels
$ws.Rows(1).RowHeight = 78

# keep the new data cells (B2:F17) empty - nothing else to do there

# move the sheet's active-cell selection, as in the authored change
$ws.Activate()
$ws.Range("D12").Select()

# restore the workbook's active sheet/selection back to the first sheet
$ws1 = $wb.Worksheets.Item("UserGuide")
$ws1.Activate()
$ws1.Range("A1").Select()
